# Insert one new data row into the daily/hourly log table on row 734,
# pushing the existing rows 734-775 down to 735-776 (dimension grows to
# A1:D776). The new row records an entry for 2026/01/29 (Thursday) that
# was missing from the export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(734).Insert()

# Column A stores dates as plain text (e.g. "2026/12/29"), not real Excel
# dates, elsewhere in this sheet. Force text interpretation so "2026/01/29"
# isn't auto-converted into a date serial, then drop back to the workbook's
# plain/default cell style so the new row matches its neighbours exactly.
$ws.Range("A734").NumberFormat = "@"
$ws.Range("A734").Value = "2026/01/29"
$ws.Range("B734").Value = "木"
$ws.Range("C734").Value = 12
$ws.Range("D734").Value = 201
$ws.Range("A734:D734").Style = "Normal"
